# Swap the weekly price-record values between row 2 and row 4.
# (Row 2 previously held the newer date/prices, row 4 the older ones;
#  the edit puts them back in chronological order.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("D", "N", "O", "P", "S")

# Read all the current values first so writes don't clobber reads.
$row2Values = @{}
$row4Values = @{}
foreach ($col in $columns) {
    $row2Values[$col] = $ws.Range("$col`2").Value2
    $row4Values[$col] = $ws.Range("$col`4").Value2
}

foreach ($col in $columns) {
    $ws.Range("$col`2").Value2 = $row4Values[$col]
    $ws.Range("$col`4").Value2 = $row2Values[$col]
}
